# Applies the "Updated cryptos list" data refresh: new Price (column D)
# and Volume(1h) (column E) values for the coinranking snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.978.44"
$ws.Range("E2").Value = "  +3.35%  "

$ws.Range("D3").Value = "2.416.02"
$ws.Range("E3").Value = "  +2.93%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.11%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.106"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.358"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.50%  "

$ws.Range("D14").Value = "2.849.59"
$ws.Range("E14").Value = "  +3.02%  "

$ws.Range("D15").Value = "59.935.94"
$ws.Range("E15").Value = "  +3.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000138"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.44%  "

$ws.Range("D17").Value = "2.406.06"
$ws.Range("E17").Value = "  +2.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "331.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "

$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.80%  "

$ws.Range("E24").Value = "  +3.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.98%  "

$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("D28").Value = "0.0₃0780"
$ws.Range("E28").Value = "  +6.00%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("E33").Value = "  +2.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  +5.21%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.415"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "312.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0962"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "

$ws.Range("E45").Value = "  +1.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.574"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.409"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0225"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
